# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.995.49'
$ws.Range("D2").Style = "Normal"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.329.91'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.496'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.17'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.29%  '
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.21'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.18%  '
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.26%  '
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.22%  '
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.695.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.361.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.93%  '
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.791'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.18%  '
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.950.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.15'
$ws.Range("D19").Style = "Normal"
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.35%  '
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.57%  '
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.84%  '
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("E27").Style = "Normal"
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -13.40%  '
$ws.Range("E28").Style = "Normal"
# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.97%  '
$ws.Range("E29").Style = "Normal"
# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.43%  '
$ws.Range("E30").Style = "Normal"
# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("E31").Style = "Normal"
# Row 32
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '138.58'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -16.42%  '
$ws.Range("E32").Style = "Normal"
# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.72%  '
$ws.Range("E33").Style = "Normal"
# Row 34
$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.76'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.12%  '
$ws.Range("E34").Style = "Normal"
# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.91%  '
$ws.Range("E35").Style = "Normal"
# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.12%  '
$ws.Range("E36").Style = "Normal"
# Row 37
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.15%  '
$ws.Range("E37").Style = "Normal"
# Row 38
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.30'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.83%  '
$ws.Range("E38").Style = "Normal"
# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.58%  '
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.28'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +23.85%  '
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.932.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.22%  '
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0279'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.66%  '
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.57%  '
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.48%  '
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.561.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.20%  '
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.68%  '
$ws.Range("E51").Style = "Normal"
